$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 116.067275999999
$ws.Range("A3").Value = 119.82653999999999
$ws.Range("A4").Value = 119.983176
$ws.Range("A5").Value = 119.82653999999999
$ws.Range("A6").Value = 119.82653999999999
$ws.Range("A7").Value = 121.079628
$ws.Range("A8").Value = 119.043359999999
$ws.Range("A9").Value = 119.82653999999999
$ws.Range("A10").Value = 120.296448
$ws.Range("A11").Value = 121.236263999999
$ws.Range("A12").Value = 119.983176
$ws.Range("A13").Value = 116.067275999999
$ws.Range("A14").Value = 121.3929
$ws.Range("A15").Value = 118.260179999999
$ws.Range("A16").Value = 120.453084
$ws.Range("A17").Value = 119.199996
$ws.Range("A18").Value = 119.82653999999999
$ws.Range("A19").Value = 119.199996
$ws.Range("A20").Value = 119.356632
$ws.Range("A21").Value = 120.766356
$ws.Range("A22").Value = 120.296448
$ws.Range("A23").Value = 119.669904
$ws.Range("A24").Value = 121.236263999999
$ws.Range("A25").Value = 119.043359999999
$ws.Range("A26").Value = 120.60972
$ws.Range("A27").Value = 120.453084
$ws.Range("A28").Value = 121.549536
$ws.Range("A29").Value = 119.356632
$ws.Range("A30").Value = 114.814188
$ws.Range("A31").Value = 122.01944399999999

$ws.Range("E6").Select()
